# Update header volume number and week-covering dates (rich text runs reconstituted as plain text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# Update weekly crime statistics table (rows 14-33)

$ws.Range("C14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 200
$ws.Range("M14").Value = -15.384615384615
$ws.Range("N14").Value = -81.967213114754
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 56
$ws.Range("K15").Value = 80.645161290322
$ws.Range("L15").Value = 43.589743589743
$ws.Range("M15").Value = 69.696969696969
$ws.Range("N15").Value = -13.846153846153
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 40
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 420
$ws.Range("J16").Value = 438
$ws.Range("K16").Value = -4.109589041095
$ws.Range("L16").Value = -11.205073995771
$ws.Range("M16").Value = -7.079646017699
$ws.Range("N16").Value = -76.910390324354
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 62
$ws.Range("H17").Value = 4.838709677419
$ws.Range("I17").Value = 705
$ws.Range("J17").Value = 688
$ws.Range("K17").Value = 2.470930232558
$ws.Range("L17").Value = -5.495978552278
$ws.Range("M17").Value = 56.319290465631
$ws.Range("N17").Value = -16.370106761565
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 238
$ws.Range("J18").Value = 274
$ws.Range("K18").Value = -13.138686131386
$ws.Range("L18").Value = -15.302491103202
$ws.Range("M18").Value = -30.205278592375
$ws.Range("N18").Value = -85.514303104077
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 77
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -10.465116279069
$ws.Range("I19").Value = 970
$ws.Range("J19").Value = 931
$ws.Range("K19").Value = 4.189044038668
$ws.Range("L19").Value = 18.581907090464
$ws.Range("M19").Value = 92.079207920792
$ws.Range("N19").Value = 39.769452449567
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 52
$ws.Range("H20").Value = -57.692307692307
$ws.Range("I20").Value = 417
$ws.Range("J20").Value = 478
$ws.Range("K20").Value = -12.76150627615
$ws.Range("L20").Value = -27.604166666666
$ws.Range("M20").Value = 93.055555555555
$ws.Range("N20").Value = -76.400679117147
$ws.Range("D21").Value = 54
$ws.Range("E21").Value = -5.555555555555
$ws.Range("F21").Value = 221
$ws.Range("G21").Value = 268
$ws.Range("H21").Value = -17.537313432835
$ws.Range("I21").Value = 2817
$ws.Range("J21").Value = 2852
$ws.Range("K21").Value = -1.227208976157
$ws.Range("L21").Value = -4.313858695652
$ws.Range("M21").Value = 40.079562406762
$ws.Range("N21").Value = -59.126523505513
$ws.Range("D22").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -100
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = 35
$ws.Range("H23").Value = -37.142857142857
$ws.Range("I23").Value = 236
$ws.Range("J23").Value = 262
$ws.Range("K23").Value = -9.923664122137
$ws.Range("L23").Value = -14.181818181818
$ws.Range("M23").Value = 18.59296482412
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -3.333333333333
$ws.Range("F24").Value = 153
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = 15.037593984962
$ws.Range("I24").Value = 1805
$ws.Range("J24").Value = 1653
$ws.Range("K24").Value = 9.19540229885
$ws.Range("L24").Value = 5.370694687682
$ws.Range("M24").Value = 37.262357414448
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 20.51282051282
$ws.Range("I25").Value = 562
$ws.Range("J25").Value = 644
$ws.Range("K25").Value = -12.732919254658
$ws.Range("L25").Value = -19.828815977175
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -11.764705882352
$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 75
$ws.Range("H26").Value = -17.333333333333
$ws.Range("I26").Value = 949
$ws.Range("J26").Value = 1039
$ws.Range("K26").Value = -8.662175168431
$ws.Range("L26").Value = -2.666666666666
$ws.Range("M26").Value = -31.031976744186
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 67
$ws.Range("K27").Value = 39.583333333333
$ws.Range("L27").Value = 17.543859649122
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 22.222222222222
$ws.Range("I28").Value = 69
$ws.Range("J28").Value = 99
$ws.Range("K28").Value = -30.30303030303
$ws.Range("L28").Value = -29.591836734693
$ws.Range("C29").Value = 1
$ws.Range("I29").Value = 35
$ws.Range("K29").Value = 34.615384615384
$ws.Range("L29").Value = -10.25641025641
$ws.Range("M29").Value = -23.91304347826
$ws.Range("N29").Value = -74.63768115942
$ws.Range("C30").Value = 1
$ws.Range("I30").Value = 30
$ws.Range("K30").Value = 30.434782608695
$ws.Range("L30").Value = -9.090909090909
$ws.Range("M30").Value = -23.076923076923
$ws.Range("N30").Value = -75.806451612903
$ws.Range("D33").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("L14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("G33").Value = 3
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = -66.666666666666

$excel.CutCopyMode = 0
